$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.127
$ws.Range("B14").Value = 5.661
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("D22").Value = -8.098000000000001
$ws.Range("B23").Value = 7.886
$ws.Range("D24").Value = -6.969000000000001
$ws.Range("B25").Value = 5.976999999999999
$ws.Range("B26").Value = 5.630000000000001
$ws.Range("D28").Value = -7.986000000000002
$ws.Range("B29").Value = 5.519
$ws.Range("D36").Value = -7.540999999999999
$ws.Range("D45").Value = -7.496
$ws.Range("D48").Value = -7.572
$ws.Range("D49").Value = -8.425000000000001
$ws.Range("D52").Value = -8.199000000000002
$ws.Range("B53").Value = 5.706
$ws.Range("D53").Value = -8.343999999999999
$ws.Range("D54").Value = -8.239000000000001
$ws.Range("B57").Value = 5.034000000000001
$ws.Range("B59").Value = 4.85
$ws.Range("B69").Value = 5.411
$ws.Range("D70").Value = -7.026999999999999
$ws.Range("B79").Value = 5.889
$ws.Range("B83").Value = 5.473999999999999
$ws.Range("D86").Value = -8.252000000000001
$ws.Range("D89").Value = -8.177000000000001
$ws.Range("B91").Value = 5.601000000000001
$ws.Range("B93").Value = 5.369000000000001
$ws.Range("D101").Value = -7.898000000000001
$ws.Range("B103").Value = 5.267
